$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 8.351000000000001
$ws.Range("C4").Value = -12.065
$ws.Range("B7").Value = 5.117
$ws.Range("D10").Value = -7.776999999999999
$ws.Range("C12").Value = -11.099
$ws.Range("D13").Value = -8.061
$ws.Range("B16").Value = 4.808999999999999
$ws.Range("C18").Value = -12.708
$ws.Range("C19").Value = -11.951
$ws.Range("C20").Value = -12.31
$ws.Range("B28").Value = 5.293000000000001
$ws.Range("B29").Value = 5.213
$ws.Range("D30").Value = -7.290000000000001
$ws.Range("C31").Value = -13.117
$ws.Range("B32").Value = 7.084999999999999
$ws.Range("B40").Value = 9.178000000000001
$ws.Range("C40").Value = -11.873
$ws.Range("D40").Value = -8.43
$ws.Range("C42").Value = -12.31
$ws.Range("D44").Value = -7.744
$ws.Range("C47").Value = -12.304
$ws.Range("C48").Value = -11.936
$ws.Range("B52").Value = 5.078999999999999
$ws.Range("B57").Value = 4.882000000000001
$ws.Range("C63").Value = -10.775
$ws.Range("C64").Value = -11.18
$ws.Range("B66").Value = 5.548
$ws.Range("C76").Value = -11.857
$ws.Range("C81").Value = -13.317
$ws.Range("C89").Value = -13.331
$ws.Range("D89").Value = -8.523999999999999
$ws.Range("D91").Value = -7.408999999999999
$ws.Range("C94").Value = -11.375
$ws.Range("B100").Value = 5.828999999999999
